# Presence database - arrival & departure functions
# Rename existing sheets, rework their "NO" availability markers, and add
# two brand-new sheets ("rotem" and "michal") at the end of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the five existing sheets (order is preserved: positions 1..5)
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "adir"   # was "michal"
$wb.Worksheets.Item(2).Name = "stav"   # was "shir"
$wb.Worksheets.Item(3).Name = "tair"   # was "emilia"
$wb.Worksheets.Item(4).Name = "yoni"   # was "emilia&shir"
$wb.Worksheets.Item(5).Name = "asaf"   # was "kakaka"

# ---------------------------------------------------------------------
# 2) adir (sheet1): remove E2 and D3, keep B2 and G3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("adir")
$ws.Range("E2").ClearContents()
$ws.Range("D3").ClearContents()

# ---------------------------------------------------------------------
# 3) stav (sheet2): move B2 -> H2, move G3 -> F3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("stav")
$ws.Range("B2").ClearContents()
$ws.Range("H2").Value = "NO"
$ws.Range("G3").ClearContents()
$ws.Range("F3").Value = "NO"

# ---------------------------------------------------------------------
# 4) tair (sheet3): remove H2; move F3 -> C3 and also add E3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("tair")
$ws.Range("H2").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("C3").Value = "NO"
$ws.Range("E3").Value = "NO"

# ---------------------------------------------------------------------
# 5) yoni (sheet4): keep B2, add G2; remove C3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("yoni")
$ws.Range("G2").Value = "NO"
$ws.Range("C3").ClearContents()

# ---------------------------------------------------------------------
# 6) asaf (sheet5): keep B2, remove G2; add C3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("asaf")
$ws.Range("G2").ClearContents()
$ws.Range("C3").Value = "NO"

# ---------------------------------------------------------------------
# 7) Add two brand new presence sheets at the end: "rotem" and "michal"
# ---------------------------------------------------------------------
function New-PresenceSheet {
    param(
        [string]$SheetName,
        [string]$MorningNoCol,
        [string]$EveningNoCol
    )

    $wbLocal = $excel.ActiveWorkbook
    $afterSheet = $wbLocal.Worksheets.Item($wbLocal.Worksheets.Count)
    $newWs = $wbLocal.Worksheets.Add($null, $afterSheet)
    $newWs.Name = $SheetName

    $newWs.Range("B1").Value = "Sunday"
    $newWs.Range("C1").Value = "Monday"
    $newWs.Range("D1").Value = "Tuesday"
    $newWs.Range("E1").Value = "Wednesday"
    $newWs.Range("F1").Value = "Thursday"
    $newWs.Range("G1").Value = "Friday"
    $newWs.Range("H1").Value = "Saturday"

    $newWs.Range("A2").Value = "Morning"
    $newWs.Range($MorningNoCol + "2").Value = "NO"

    $newWs.Range("A3").Value = "Evening"
    $newWs.Range($EveningNoCol + "3").Value = "NO"
}

New-PresenceSheet "rotem"  "E" "H"
New-PresenceSheet "michal" "B" "G"
